$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# This workbook is a rolling "trailing quarters" database extract:
# every update drops the oldest reporting period (column D) and
# appends the newest one (new column M), shifting everything else
# one column to the left - history is preserved from E..M -> D..L.
# ------------------------------------------------------------------

# Step 1: drop the oldest period column (D = "6 ماهه منتهی به 1399/06")
# Deleting shifts E:M left into D:L automatically (values, styles,
# shared-string usage, and the width groupings that track the
# "annual report" columns all move with it).
$ws.Columns.Item(4).Delete()

# Step 2: create the new trailing column M by cloning the formatting
# of the (now) last existing column L.
$ws.Range("L1:L28").Copy($ws.Range("M1:M28"))
$ws.Columns.Item(13).ColumnWidth = 28.1666666666667

# Step 3: header row - new period label / publish date for column M
$ws.Range("M8").Value = "12 ماهه منتهی به 1401/12"
$ws.Range("M9").Value = "1402-02-27 (3)"

# Step 4: new financial figures for the "12 ماهه منتهی به 1401/12" column
$ws.Range("M11").Value = 21157343
$ws.Range("M12").Value = -8863829
$ws.Range("M13").Value = 12293514
$ws.Range("M14").Value = -1806898
$ws.Range("M15").Value = 0
$ws.Range("M16").Value = 379865
$ws.Range("M17").Value = 10866481
$ws.Range("M18").Value = -55460
$ws.Range("M19").Value = 1250178
$ws.Range("M20").Value = 12061199
$ws.Range("M21").Value = -1766521
$ws.Range("M22").Value = 10294678
$ws.Range("M23").Value = 0
$ws.Range("M24").Value = 10294678
$ws.Range("M25").Value = 15838
$ws.Range("M26").Value = 650000
$ws.Range("M27").Value = 15838

Write-Output "edit complete"
